# Update the theta_se (row 4) and lambda_se (row 6) standard-error values
# on the sole worksheet to reflect the refreshed replication results.
# Values are written column by column (theta_se then lambda_se for each
# column) so new shared-string entries are created in the same order as
# the source replication data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..J correspond to hours 20,21,22,23,0,1,2,3,4
$ws.Range("B4").Value = "(2.22)"
$ws.Range("B6").Value = "(1.48)"

$ws.Range("C4").Value = "(1.37)"
$ws.Range("C6").Value = "(0.98)"

$ws.Range("D4").Value = "(5.16)"
$ws.Range("D6").Value = "(2.34)"

$ws.Range("E4").Value = "(0.73)"
$ws.Range("E6").Value = "(0.09)"

$ws.Range("F4").Value = "(3.41)"
$ws.Range("F6").Value = "(1.64)"

$ws.Range("G4").Value = "(0.43)"
$ws.Range("G6").Value = "(0.18)"

$ws.Range("H4").Value = "(0.15)"
$ws.Range("H6").Value = "(0.28)"

$ws.Range("I4").Value = "(0.81)"
$ws.Range("I6").Value = "(0.77)"

$ws.Range("J4").Value = "(0.15)"
$ws.Range("J6").Value = "(1.12)"
